$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update country names (column A) for rows whose rank changed due to re-sort by "Casos totales"
$ws.Range("A30").Value = 'Bielorrusia'
$ws.Range("A31").Value = 'Catar'
$ws.Range("A73").Value = 'Azerbaiyan'
$ws.Range("A74").Value = 'Camerun'
$ws.Range("A118").Value = 'Tayikistan'
$ws.Range("A119").Value = 'Gabon'
$ws.Range("A120").Value = 'Jordania'
$ws.Range("A121").Value = 'Malta'
$ws.Range("A122").Value = 'Jamaica'
$ws.Range("A123").Value = 'Tanzania'
$ws.Range("A124").Value = 'Paraguay'
$ws.Range("A147").Value = 'Zambia'
$ws.Range("A148").Value = 'Suazilandia'
$ws.Range("A165").Value = 'Mozambique'
$ws.Range("A166").Value = 'Liechtenstein'
$ws.Range("A167").Value = 'Barbados'
$ws.Range("A181").Value = 'Yemen'
$ws.Range("A182").Value = 'Zimbabue'
$ws.Range("A183").Value = 'Guam'
$ws.Range("A192").Value = 'Belice'
$ws.Range("A193").Value = 'Nueva Caledonia'

# Update updated case-count statistics (columns B:H) for affected rows
$ws.Range("B4").Value = 1295101
$ws.Range("C4").Value = 2478
$ws.Range("D4").Value = 217292
$ws.Range("E4").Value = 1000750
$ws.Range("F4").Value = 16992
$ws.Range("G4").Value = 131
$ws.Range("H4").Value = 77059

$ws.Range("F19").Value = 564

$ws.Range("B30").Value = 21101
$ws.Range("C30").Value = 933
$ws.Range("D30").Value = 5484
$ws.Range("E30").Value = 15496
$ws.Range("F30").Value = 92
$ws.Range("G30").Value = 5
$ws.Range("H30").Value = 121

$ws.Range("B31").Value = 20201
$ws.Range("C31").Value = 1311
$ws.Range("D31").Value = 2370
$ws.Range("E31").Value = 17819
$ws.Range("F31").Value = 72
$ws.Range("H31").Value = 12

$ws.Range("F44").Value = 45

$ws.Range("B46").Value = 9376
$ws.Range("C46").Value = 281
$ws.Range("D46").Value = 2286
$ws.Range("E46").Value = 6710
$ws.Range("F46").Value = 134
$ws.Range("G46").Value = 7
$ws.Range("H46").Value = 380

$ws.Range("D55").Value = 4000
$ws.Range("E55").Value = 1478

$ws.Range("B60").Value = 4728
$ws.Range("C60").Value = 123
$ws.Range("D60").Value = 1826
$ws.Range("E60").Value = 2752
$ws.Range("G60").Value = 5
$ws.Range("H60").Value = 150

$ws.Range("D72").Value = 1775
$ws.Range("E72").Value = 539

$ws.Range("B73").Value = 2279
$ws.Range("C73").Value = 75
$ws.Range("D73").Value = 1576
$ws.Range("E73").Value = 675
$ws.Range("F73").Value = 18
$ws.Range("H73").Value = 28

$ws.Range("B74").Value = 2267
$ws.Range("D74").Value = 1002
$ws.Range("E74").Value = 1157
$ws.Range("F74").Value = 12
$ws.Range("H74").Value = 108

$ws.Range("B79").Value = 1872
$ws.Range("C79").Value = 43
$ws.Range("E79").Value = 1385
$ws.Range("F79").Value = 49
$ws.Range("G79").Value = 2
$ws.Range("H79").Value = 86

$ws.Range("B81").Value = 1741
$ws.Range("C81").Value = 12
$ws.Range("D81").Value = 1078
$ws.Range("E81").Value = 589
$ws.Range("G81").Value = 1
$ws.Range("H81").Value = 74

$ws.Range("B89").Value = 1450
$ws.Range("C89").Value = 1
$ws.Range("D89").Value = 252
$ws.Range("E89").Value = 1098
$ws.Range("F89").Value = 12
$ws.Range("G89").Value = 1
$ws.Range("H89").Value = 100

$ws.Range("B118").Value = 522
$ws.Range("C118").Value = 61
$ws.Range("D118").Value = 0
$ws.Range("E118").Value = 510
$ws.Range("F118").Value = 0
$ws.Range("H118").Value = 12

$ws.Range("B119").Value = 504
$ws.Range("D119").Value = 110
$ws.Range("E119").Value = 386
$ws.Range("F119").Value = 1
$ws.Range("H119").Value = 8

$ws.Range("B120").Value = 494
$ws.Range("C120").Value = 0
$ws.Range("D120").Value = 381
$ws.Range("E120").Value = 104
$ws.Range("F120").Value = 5
$ws.Range("H120").Value = 9

$ws.Range("B121").Value = 489
$ws.Range("C121").Value = 3
$ws.Range("D121").Value = 419
$ws.Range("E121").Value = 65
$ws.Range("H121").Value = 5

$ws.Range("B122").Value = 488
$ws.Range("C122").Value = 10
$ws.Range("D122").Value = 58
$ws.Range("E122").Value = 421
$ws.Range("F122").Value = 0
$ws.Range("H122").Value = 9

$ws.Range("B123").Value = 480
$ws.Range("D123").Value = 167
$ws.Range("E123").Value = 297
$ws.Range("F123").Value = 7
$ws.Range("H123").Value = 16

$ws.Range("B124").Value = 462
$ws.Range("D124").Value = 148
$ws.Range("E124").Value = 304
$ws.Range("F124").Value = 9
$ws.Range("H124").Value = 10

$ws.Range("D129").Value = 226
$ws.Range("E129").Value = 147

$ws.Range("B139").Value = 230
$ws.Range("C139").Value = 12
$ws.Range("D139").Value = 44
$ws.Range("E139").Value = 184

$ws.Range("D146").Value = 67
$ws.Range("E146").Value = 103

$ws.Range("B147").Value = 167
$ws.Range("C147").Value = 14
$ws.Range("D147").Value = 111
$ws.Range("E147").Value = 52
$ws.Range("F147").Value = 1
$ws.Range("H147").Value = 4

$ws.Range("D148").Value = 12
$ws.Range("E148").Value = 139
$ws.Range("F148").Value = 0
$ws.Range("H148").Value = 2

$ws.Range("C165").Value = 1
$ws.Range("D165").Value = 27
$ws.Range("E165").Value = 55
$ws.Range("H165").Value = 0

$ws.Range("D166").Value = 55
$ws.Range("E166").Value = 26
$ws.Range("F166").Value = 0
$ws.Range("H166").Value = 1

$ws.Range("B167").Value = 82
$ws.Range("D167").Value = 53
$ws.Range("E167").Value = 22
$ws.Range("F167").Value = 4
$ws.Range("H167").Value = 7

$ws.Range("C181").Value = 9
$ws.Range("D181").Value = 1
$ws.Range("E181").Value = 26
$ws.Range("G181").Value = 2
$ws.Range("H181").Value = 7

$ws.Range("B182").Value = 34
$ws.Range("D182").Value = 5
$ws.Range("E182").Value = 25
$ws.Range("H182").Value = 4

$ws.Range("B183").Value = 32
$ws.Range("D183").Value = 0
$ws.Range("E183").Value = 31
$ws.Range("H183").Value = 1

$ws.Range("D192").Value = 16
$ws.Range("H192").Value = 2

$ws.Range("D193").Value = 18

# Update the "last updated" timestamp string
$ws.Range("A1").Value = "Datos actualizados a 8 de Mayo de 2020 a las 17:04"

Write-Host "Edit applied"